$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the Fecha (D) and Volumen (J) values between row 2 and row 5
$ws.Range("D2").Value = 44714
$ws.Range("J2").Value = 80

$ws.Range("D5").Value = 44792
$ws.Range("J5").Value = 160
